# Apply the refreshed crypto price/volume snapshot to Sheet1.
# Values that look like plain numbers (e.g. "563.21") are written with a
# leading apostrophe so Excel keeps them as text, matching the workbook's
# existing convention of storing the Price column as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.186.51'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.477.22'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('D5').Value = '''563.21'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '''163.21'
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').Value = '2.474.74'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').Value = '''0.152'
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('E12').Value = '  -1.91%  '
$ws.Range('E13').Value = '  +1.86%  '
$ws.Range('D14').Value = '69.056.58'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Value = '''23.74'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = '2.478.18'
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').Value = '''10.81'
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('D20').Value = '''340.04'
$ws.Range('E20').Value = '  -2.39%  '
$ws.Range('D21').Value = '''7.04'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = '''1.91'
$ws.Range('E23').Value = '  +3.68%  '
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').Value = '''67.44'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('D27').Value = '2.604.56'
$ws.Range('E27').Value = '  +1.44%  '
$ws.Range('D28').Value = '''8.28'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').Value = '''0.997'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').Value = '0.0₃0825'
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('D31').Value = '''7.22'
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('D32').Value = '''435.73'
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('E35').Value = '  -2.14%  '
$ws.Range('D36').Value = '''157.52'
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Value = '''0.302'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('D42').Value = '''4.46'
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''1.48'
$ws.Range('E43').Value = '  -2.71%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').Value = '''1.08'
$ws.Range('E44').Value = '  +1.79%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''2.09'
$ws.Range('E45').Value = '  +1.36%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''133.75'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '''3.36'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.0718'
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '''0.487'
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.565'
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.0915'
$ws.Range('E51').Value = '  +0.22%  '
